$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# A new item ("شاش فازلين 10*10 سم") was sold and needs to be inserted into
# the items table. It takes the place of row 9 (between "جونتي عمال" and
# "معطر جو FRIDA "), and the item that used to be row 9 ("معطر جو FRIDA ")
# is pushed down to a new row 10. The totals row and the footer row (old
# rows 10 and 11) shift down to rows 11 and 12.
# ---------------------------------------------------------------------------

# Insert a new blank row at position 10, pushing the totals/footer rows down
$ws.Rows(10).Insert()

# Copy the formatting (styles, number formats, etc.) of row 9 into the new
# row 10 so it looks just like the other item rows, then fix up the row
# heights to match the final layout.
$ws.Range("A9:Q9").Copy()
$ws.Range("A10:Q10").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows(10).RowHeight = 24.75
$ws.Rows(11).RowHeight = 25.5

# Recreate the merged cells for row 10 (same merge pattern as row 9).
$ws.Range("A10:B10").Merge()
$ws.Range("C10:G10").Merge()
$ws.Range("H10:K10").Merge()
$ws.Range("L10:M10").Merge()
$ws.Range("N10:O10").Merge()

# Row 10 now holds the item that used to be in row 9: معطر جو FRIDA
$ws.Range("A10").Value2 = 4
$ws.Range("C10").Value2 = "معطر جو FRIDA "
$ws.Range("H10").Value2 = "7:0"
$fmt = $ws.Range("L10").NumberFormat
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value2 = "0"
$ws.Range("L10").NumberFormat = $fmt
$ws.Range("N10").Value2 = "65.00"
$fmt = $ws.Range("P10").NumberFormat
$ws.Range("P10").NumberFormat = "@"
$ws.Range("P10").Value2 = "65.0000"
$ws.Range("P10").NumberFormat = $fmt
$ws.Range("Q10").Value2 = "1:0"

# Row 9 becomes the newly added item: شاش فازلين 10*10 سم
$ws.Range("C9").Value2 = "شاش فازلين 10*10 سم"
$ws.Range("H9").Value2 = "7:0"
$fmt = $ws.Range("L9").NumberFormat
$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value2 = "0"
$ws.Range("L9").NumberFormat = $fmt
$ws.Range("N9").Value2 = "7.00"
$fmt = $ws.Range("P9").NumberFormat
$ws.Range("P9").NumberFormat = "@"
$ws.Range("P9").Value2 = "7.0000"
$ws.Range("P9").NumberFormat = $fmt
$ws.Range("Q9").Value2 = "1:0"

# Update the totals row (now row 11): add the new item's price (7.00) to the
# previous total of 184.2 -> 191.2
$ws.Range("P11").Value2 = 191.2

# Update the generated-on timestamp shown in the footer (now row 12)
$ws.Range("A12").Value2 = "Sunday, 21 September, 2025 11:17 AM"
